$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 214, shifting existing rows
# 214-287 down to 216-289.
$ws.Rows.Item(214).Insert()
$ws.Rows.Item(214).Insert()

# --- New row 214 ---
$ws.Cells.Item(214,1).Value = 5
$ws.Cells.Item(214,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(214,3).Value = "Maule"
$ws.Cells.Item(214,4).Value = 44524
$ws.Cells.Item(214,5).Value = 7
$ws.Cells.Item(214,6).Value = "Fruta"
$ws.Cells.Item(214,7).Value = 100109
$ws.Cells.Item(214,8).Value = "Uva"
$ws.Cells.Item(214,9).Value = 100109001
$ws.Cells.Item(214,10).Value = "Uva"
$ws.Cells.Item(214,11).Value = "Flame Seedless"
$ws.Cells.Item(214,12).Value = "Primera"
$ws.Cells.Item(214,13).Value = 280
$ws.Cells.Item(214,14).Value = 22000
$ws.Cells.Item(214,15).Value = 22000
$ws.Cells.Item(214,16).Value = 22000
$ws.Cells.Item(214,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(214,18).Value = "Provincia del Elquí"
$ws.Cells.Item(214,19).Value = 2200
$ws.Cells.Item(214,20).Value = 10

# --- New row 215 ---
$ws.Cells.Item(215,1).Value = 5
$ws.Cells.Item(215,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(215,3).Value = "Maule"
$ws.Cells.Item(215,4).Value = 44524
$ws.Cells.Item(215,5).Value = 7
$ws.Cells.Item(215,6).Value = "Fruta"
$ws.Cells.Item(215,7).Value = 100109
$ws.Cells.Item(215,8).Value = "Uva"
$ws.Cells.Item(215,9).Value = 100109001
$ws.Cells.Item(215,10).Value = "Uva"
$ws.Cells.Item(215,11).Value = "Red Globe"
$ws.Cells.Item(215,12).Value = "Primera"
$ws.Cells.Item(215,13).Value = 150
$ws.Cells.Item(215,14).Value = 25000
$ws.Cells.Item(215,15).Value = 25000
$ws.Cells.Item(215,16).Value = 25000
$ws.Cells.Item(215,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(215,18).Value = "Provincia del Elquí"
$ws.Cells.Item(215,19).Value = 2500
$ws.Cells.Item(215,20).Value = 10
